# Apply a cyclic rotation of the data rows 4-7 in the "Artfynd" worksheet.
# Row 7's data moves to row 4, row 4's data moves to row 5,
# row 5's data moves to row 6, and row 6's data moves to row 7.
# Columns C, D, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY
# (location / metadata columns) are identical across these rows and are
# left untouched. Column I is empty for all of these rows too.
#
# Because row 7 previously used column M (Aktivitet) instead of column J
# (Enhet), and the other rows use column J instead of M, we must clear
# out the cell that is no longer used for each destination row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for rows 4-7 so that the rotation source
# data does not get clobbered while we write the new values.
$rowsBefore = @{}
foreach ($r in 4..7) {
    $rowsBefore[$r] = @{
        A = $ws.Cells.Item($r, 1).Value2
        B = $ws.Cells.Item($r, 2).Value2
        E = $ws.Cells.Item($r, 5).Value2
        F = $ws.Cells.Item($r, 6).Value2
        G = $ws.Cells.Item($r, 7).Value2
        H = $ws.Cells.Item($r, 8).Value2
        J = $ws.Cells.Item($r, 10).Value2
        M = $ws.Cells.Item($r, 13).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
    }
}

# Mapping of destination row -> source row (cyclic rotation).
$sourceForDest = @{ 4 = 7; 5 = 4; 6 = 5; 7 = 6 }

foreach ($dest in 4..7) {
    $src = $sourceForDest[$dest]
    $data = $rowsBefore[$src]

    $ws.Cells.Item($dest, 1).Value = $data.A   # A - Id
    $ws.Cells.Item($dest, 2).Value = $data.B   # B - Taxonsorteringsordning
    $ws.Cells.Item($dest, 5).Value = $data.E   # E - TaxonId
    $ws.Cells.Item($dest, 6).Value = $data.F   # F - Artnamn
    $ws.Cells.Item($dest, 7).Value = $data.G   # G - Vetenskapligt namn
    $ws.Cells.Item($dest, 8).Value = $data.H   # H - Auktor
    $ws.Cells.Item($dest, 17).Value = $data.Q  # Q - Ost
    $ws.Cells.Item($dest, 18).Value = $data.R  # R - Nord

    # J (Enhet) and M (Aktivitet) are mutually exclusive in these rows:
    # only one of them is populated per row. Clear both, then set
    # whichever one the source row actually used.
    $ws.Cells.Item($dest, 10).Value = ""
    $ws.Cells.Item($dest, 13).Value = ""

    if ($data.J -ne $null -and $data.J -ne "") {
        $ws.Cells.Item($dest, 10).Value = $data.J
    }
    if ($data.M -ne $null -and $data.M -ne "") {
        $ws.Cells.Item($dest, 13).Value = $data.M
    }
}
